$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking Price cells first, so Excel stores the new values
# as text (matching the original "Price" column cell type) instead of auto-converting them
# to floating point numbers.
$textCells = @("D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D18","D19","D20","D21","D22","D24","D25","D26","D27","D28","D29","D30","D31","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated Price / Volume(1h) figures (and the two Coin/Link row swaps)
$ws.Range("D2").Value = "23.421.91"
$ws.Range("E2").Value = "  +0.98%  "
$ws.Range("D3").Value = "1.638.33"
$ws.Range("E3").Value = "  +2.21%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").Value = "304.79"
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("D7").Value = "0.3718"
$ws.Range("E7").Value = "  -1.69%  "
$ws.Range("D8").Value = "51.86"
$ws.Range("E8").Value = "  -0.33%  "
$ws.Range("D9").Value = "0.3612"
$ws.Range("E9").Value = "  -0.40%  "
$ws.Range("D10").Value = "1.249"
$ws.Range("E10").Value = "  -1.75%  "
$ws.Range("D11").Value = "0.08108"
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").Value = "22.72"
$ws.Range("E13").Value = "  -0.55%  "
$ws.Range("D14").Value = "6.569"
$ws.Range("E14").Value = "  -0.55%  "
$ws.Range("D15").Value = "0.00001264"
$ws.Range("E15").Value = "  +1.43%  "
$ws.Range("D16").Value = "7.258"
$ws.Range("E16").Value = "  -2.17%  "
$ws.Range("D17").Value = "1.633.93"
$ws.Range("E17").Value = "  +1.91%  "
$ws.Range("D18").Value = "94.09"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("D19").Value = "0.06886"
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("D20").Value = "18.06"
$ws.Range("E20").Value = "  -0.32%  "
$ws.Range("D21").Value = "6.490"
$ws.Range("E21").Value = "  -0.89%  "
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "23.433.25"
$ws.Range("E23").Value = "  +1.06%  "
$ws.Range("D24").Value = "12.69"
$ws.Range("E24").Value = "  -2.23%  "
$ws.Range("D25").Value = "2.409"
$ws.Range("E25").Value = "  +0.68%  "
$ws.Range("D26").Value = "3.032"
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").Value = "21.16"
$ws.Range("E27").Value = "  -0.36%  "
$ws.Range("D28").Value = "151.57"
$ws.Range("D29").Value = "5.309"
$ws.Range("E29").Value = "  +1.19%  "
$ws.Range("D30").Value = "135.16"
$ws.Range("E30").Value = "  +0.96%  "
$ws.Range("D31").Value = "2.278"
$ws.Range("E31").Value = "  -3.66%  "
$ws.Range("D32").Value = "1.814.92"
$ws.Range("E32").Value = "  +1.93%  "
$ws.Range("D33").Value = "6.727"
$ws.Range("E33").Value = "  -0.53%  "
$ws.Range("D34").Value = "0.9493"
$ws.Range("E34").Value = "  -1.88%  "
$ws.Range("D35").Value = "0.02799"
$ws.Range("E35").Value = "  +2.45%  "
$ws.Range("D36").Value = "10.25"
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("B37").Value = "Algorand"
$ws.Range("C37").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D37").Value = "0.2511"
$ws.Range("E37").Value = "  -0.58%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "0.07205"
$ws.Range("E38").Value = "  -4.15%  "
$ws.Range("D39").Value = "0.08755"
$ws.Range("E39").Value = "  -0.44%  "
$ws.Range("D40").Value = "6.036"
$ws.Range("E40").Value = "  -1.05%  "
$ws.Range("D41").Value = "1.370"
$ws.Range("E41").Value = "  +0.35%  "
$ws.Range("D42").Value = "0.7021"
$ws.Range("E42").Value = "  -1.32%  "
$ws.Range("D43").Value = "12.39"
$ws.Range("E43").Value = "  -2.06%  "
$ws.Range("D44").Value = "16.03"
$ws.Range("E44").Value = "  +2.27%  "
$ws.Range("D45").Value = "0.6472"
$ws.Range("E45").Value = "  -1.34%  "
$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  +0.14%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "2.313"
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("D48").Value = "3.998"
$ws.Range("E48").Value = "  -0.56%  "
$ws.Range("D49").Value = "0.07963"
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("D50").Value = "127.99"
$ws.Range("E50").Value = "  -3.32%  "
$ws.Range("D51").Value = "1.197"
$ws.Range("E51").Value = "  -0.81%  "

# Restore default ("Normal") styling on the cells we temporarily forced to text format,
# so we do not leave behind an unused/differing cell style vs. the original workbook.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
